# "added Logic in Excel" -- populate the CRD column on the IFS worksheet
# with real text values ("CRD" / "No CRD" / "Nothing") instead of the
# placeholder TRUE/FALSE booleans, and refresh the view selections that
# were left on a few sheets after the edit.

$wb = $excel.ActiveWorkbook

# --- IFS sheet: column C (CRD) -------------------------------------------
$wsIFS = $wb.Worksheets.Item("IFS")

$wsIFS.Range("C3").Value  = "CRD"
$wsIFS.Range("C4").Value  = "CRD"
$wsIFS.Range("C5").Value  = "No CRD"
$wsIFS.Range("C6").Value  = "No CRD"
$wsIFS.Range("C7").Value  = "CRD"
$wsIFS.Range("C8").Value  = "No CRD"
$wsIFS.Range("C9").Value  = "No CRD"
$wsIFS.Range("C10").Value = "No CRD"
$wsIFS.Range("C11").Value = "CRD"
$wsIFS.Range("C12").Value = "Nothing"

# leave the cursor where the author left it
$wsIFS.Range("H11").Select()

# --- SUMIF, COUNTIF, SUMIFS sheet: just a selection move ------------------
$wsSUMIF = $wb.Worksheets.Item("SUMIF, COUNTIF. SUMIFS")
$wsSUMIF.Range("I19").Select()

# --- Example applications sheet: stays the active tab/selection -----------
$wsEx = $wb.Worksheets.Item("Example applications")
$wsEx.Range("E44").Select()
